$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove stray empty cells D1/E1 left over from header row
$ws.Range("D1:E1").ClearContents()

# Delete trailing rows 224-245 (removed articles)
$ws.Range("A224:C245").EntireRow.Delete()

# Trim leading/trailing whitespace from specific title cells
$ws.Range("A32").Value = 'آزادکشمیر میں وزارت عظمیٰ کیلئے پارٹی میں اختلافات، پی ٹی آئی نام فائنل نہ کرسکی'
$ws.Range("A41").Value = 'ن لیگ کی وزیراعظم آزاد کشمیر کیلئے امیدوار لانے پر مشاورت'
$ws.Range("A42").Value = 'آج پاکستان میں جو حالات ہیں اس کی وجہ کیا ہے؟ حنیف عباسی'
$ws.Range("A45").Value = 'وزیراعظم شہباز شریف چینی کی اسمگلنگ اور ناجائز منافع خوری پر برہم'
$ws.Range("A53").Value = 'ڈانس اور میں ایک دوسرے کے دشمن ہیں: فواد خان'
$ws.Range("A75").Value = 'اسٹیٹ بینک حکومت کے اکاؤنٹ سے کیسے پیسے لے سکتا ہے؟، مفتاح اسماعیل'
$ws.Range("A76").Value = 'پاکستان مزید مالی یقین دہانی کرائے، آئی ایم ایف'
$ws.Range("A77").Value = '’’عمران خان نے مجھے فون کرکے کہا تھا کہ میں امریکہ کا مخالف نہیں ہوں‘‘'
$ws.Range("A79").Value = 'شرائط پوری، IMF کے پاس معاہدہ نہ کرنیکا کوئی بہانہ نہیں، وزیراعظم'
$ws.Range("A178").Value = 'کاؤنٹی چیمپئن شپ میں محمد عباس کی پھر شاندار پرفارمنس'
$ws.Range("A203").Value = '’میں اسے دیوار پر دے مارنا چاہتا ہوں‘: والدین کے ہاتھوں 10 ماہ کے بچے کا قتل جس نے جج کو بھی رلا دیا'
